$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.121.99"
$ws.Range("E2").Value = "  +0.89%  "

# Row 3
$ws.Range("D3").Value = "2.660.02"
$ws.Range("E3").Value = "  +1.65%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.90"
$ws.Range("E5").Value = "  +4.43%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.60"
$ws.Range("E6").Value = "  +1.30%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +0.80%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.59"
$ws.Range("E9").Value = "  -1.81%  "

# Row 10
$ws.Range("E10").Value = "  +5.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.353"
$ws.Range("E11").Value = "  +1.88%  "

# Row 12
$ws.Range("E12").Value = "  +0.01%  "

# Row 13
$ws.Range("D13").Value = "3.116.30"
$ws.Range("E13").Value = "  +1.35%  "

# Row 14
$ws.Range("D14").Value = "61.075.59"
$ws.Range("E14").Value = "  +0.92%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.13"
$ws.Range("E15").Value = "  +2.30%  "

# Row 16
$ws.Range("E16").Value = "  +2.58%  "

# Row 17
$ws.Range("D17").Value = "2.665.51"
$ws.Range("E17").Value = "  +1.38%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.78"
$ws.Range("E18").Value = "  +0.25%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "356.02"
$ws.Range("E19").Value = "  +1.28%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.70"
$ws.Range("E20").Value = "  +0.84%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.27"
$ws.Range("E21").Value = "  +1.59%  "

# Row 22
$ws.Range("E22").Value = "  -0.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.70"
$ws.Range("E23").Value = "  +1.82%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.433"
$ws.Range("E24").Value = "  +2.14%  "

# Row 25
$ws.Range("E25").Value = "  +1.47%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.22%  "

# Row 27
$ws.Range("D27").Value = "0.0₃0862"
$ws.Range("E27").Value = "  +2.29%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.41"
$ws.Range("E28").Value = "  +0.47%  "

# Row 29
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("E30").Value = "  +6.96%  "

# Row 31
$ws.Range("E31").Value = "  +4.44%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.61"
$ws.Range("E32").Value = "  +0.80%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.04"
$ws.Range("E33").Value = "  -0.39%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.16"
$ws.Range("E34").Value = "  +3.85%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.21"
$ws.Range("E35").Value = "  +1.45%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.924"
$ws.Range("E36").Value = "  +8.92%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.892"
$ws.Range("E37").Value = "  +0.49%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "309.55"
$ws.Range("E38").Value = "  +5.07%  "

# Row 39
$ws.Range("E39").Value = "  +0.84%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.83"
$ws.Range("E40").Value = "  +1.51%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.649"
$ws.Range("E41").Value = "  +3.87%  "

# Row 42
$ws.Range("E42").Value = "  +1.80%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0566"
$ws.Range("E43").Value = "  +1.89%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.08"
$ws.Range("E44").Value = "  +1.09%  "

# Row 45
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  +0.05%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.06"
$ws.Range("E46").Value = "  +3.16%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0240"
$ws.Range("E47").Value = "  +2.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.25"
$ws.Range("E48").Value = "  +8.11%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.36"
$ws.Range("E49").Value = "  +0.44%  "

# Row 50
$ws.Range("D50").Value = "1.999.11"
$ws.Range("E50").Value = "  -0.24%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.85"
$ws.Range("E51").Value = "  +2.52%  "
